# feat: add 2022-Q3 data
#
# The "总计" (overview) sheet gets a fresh top row for 2022-Q3 (the older
# quarters all slide down one row, keeping their own values), and a brand
# new worksheet "2022-Q3" is inserted right after "总计" holding the
# per-fund breakdown for that quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Overview sheet ("总计") - rewrite rows 2..7 in their final layout.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Give the new bottom row (row 7) the same look (bold/centered/bordered)
# as the other index cells in column A before filling it with data.
$summary.Range("A2").Copy()
$summary.Range("A7").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 7
$summary.Range("D2").Value = 1.99

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 6
$summary.Range("D3").Value = 2.35

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2021-Q4"
$summary.Range("C4").Value = 13
$summary.Range("D4").Value = 1.9

$summary.Range("A5").Value = 3
$summary.Range("B5").Value = "2021-Q2"
$summary.Range("C5").Value = 2
$summary.Range("D5").Value = 1.65

$summary.Range("A6").Value = 4
$summary.Range("B6").Value = "2021-Q1"
$summary.Range("C6").Value = 4
$summary.Range("D6").Value = 0.29

$summary.Range("A7").Value = 5
$summary.Range("B7").Value = "2020-Q4"
$summary.Range("C7").Value = 2
$summary.Range("D7").Value = 1.16

# ---------------------------------------------------------------------
# 2) New worksheet "2022-Q3", inserted right after "总计".
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $summary)
$q3.Name = "2022-Q3"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$cols = @("B", "C", "D", "E", "F", "G", "H")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $q3.Range($cols[$i] + "1")
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$rows = @(
    @(0, "360006", "光大保德信新增长混合", "17.68", "83.30", "4.51", "0.7974", 2),
    @(1, "011104", "光大保德信智能汽车主题股票", "10.24", "90.84", "5.31", "0.5437", 3),
    @(2, "001740", "光大保德信中国制造2025灵活配置混合", "10.24", "91.17", "3.79", "0.3881", 5),
    @(3, "010676", "光大保德信新机遇混合", "2.85", "84.08", "4.73", "0.1348", 8),
    @(4, "008313", "光大保德信研究精选混合", "2.26", "83.37", "4.82", "0.1089", 4),
    @(5, "000531", "东吴阿尔法灵活配置混合A", "0.28", "90.94", "5.01", "0.0140", 8),
    @(6, "014581", "东吴阿尔法灵活配置混合C", "0.03", "90.94", "5.01", "0.0015", 8)
)

$r = 2
foreach ($row in $rows) {
    $aCell = $q3.Range("A" + $r)
    $aCell.Value = $row[0]
    $aCell.Font.Bold = $true
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160
    $aCell.Borders.LineStyle = 1

    $bCell = $q3.Range("B" + $r)
    $bCell.NumberFormat = "@"
    $bCell.Value = $row[1]

    $q3.Range("C" + $r).Value = $row[2]

    $dCell = $q3.Range("D" + $r)
    $dCell.NumberFormat = "@"
    $dCell.Value = $row[3]

    $eCell = $q3.Range("E" + $r)
    $eCell.NumberFormat = "@"
    $eCell.Value = $row[4]

    $fCell = $q3.Range("F" + $r)
    $fCell.NumberFormat = "@"
    $fCell.Value = $row[5]

    $gCell = $q3.Range("G" + $r)
    $gCell.NumberFormat = "@"
    $gCell.Value = $row[6]

    $q3.Range("H" + $r).Value = $row[7]

    $r = $r + 1
}

# Restore the tab selection to the last sheet (2020-Q4), matching the
# original workbook (Worksheets.Add activates the freshly inserted sheet,
# which would otherwise leave "2022-Q3" selected).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()

